# Refresh the "cryptos" price list: update Price (D) / Volume(1h) (E) figures,
# and re-sync two rows whose coins were re-ranked and swapped places
# (Monero <-> PolygonEcosystemToken at rows 37/38, InjectiveProtocol <-> Mantle
# at rows 48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value even when it looks like a
# plain number (e.g. "573.19"), matching the source inlineStr cell type
# instead of letting Excel auto-convert it to a Number.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "62.602.82"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "2.457.85"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue "D5" "573.19"
$ws.Range("E5").Value = "  -0.64%  "

Set-TextValue "D6" "146.80"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.50%  "

$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("E12").Value = "  -0.57%  "

Set-TextValue "D13" "29.00"
$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("E14").Value = "  -2.21%  "

$ws.Range("D15").Value = "2.903.66"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").Value = "62.566.26"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "2.458.90"
$ws.Range("E17").Value = "  -0.09%  "

Set-TextValue "D18" "7.90"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("E19").Value = "  -1.59%  "

Set-TextValue "D20" "324.36"
$ws.Range("E20").Value = "  -1.80%  "

$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("E22").Value = "  +2.00%  "

Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.08%  "

Set-TextValue "D24" "10.00"
$ws.Range("E24").Value = "  +15.78%  "

$ws.Range("E25").Value = "  -1.78%  "

Set-TextValue "D26" "637.62"
$ws.Range("E26").Value = "  -3.03%  "

$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("E28").Value = "  -4.29%  "

Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  -19.59%  "

$ws.Range("E30").Value = "  -1.26%  "

Set-TextValue "D31" "7.90"
$ws.Range("E31").Value = "  -3.98%  "

$ws.Range("E33").Value = "  -4.04%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  +1.74%  "

$ws.Range("E36").Value = "  -1.50%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "151.61"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D38" "0.367"
$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("E39").Value = "  -1.47%  "

Set-TextValue "D40" "5.31"
$ws.Range("E40").Value = "  -4.43%  "

Set-TextValue "D41" "2.71"
$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("D44").Value = "0.0₆0303"
$ws.Range("E44").Value = "  -22.74%  "

Set-TextValue "D45" "152.81"
$ws.Range("E45").Value = "  +4.05%  "

Set-TextValue "D46" "15.29"
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "20.23"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "0.605"
$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("E51").Value = "  -1.81%  "
